$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of row 15 (A:E) down into the two brand-new rows (16 and 17)
# before writing any values, so the new "index" cells in column A keep the bold/
# border style used by the rest of the table.
$ws.Range("A15:E15").Copy() | Out-Null
$ws.Range("A16:E17").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Column A (row index) for the two new rows ---
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(17, 1).Value = 15

# --- Column B (row label) for every row in the table, rewritten in full so the
#     shared-string table ends up with line7/line8 in the correct sequence ---
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(10, 2).Value = "extr1"
$ws.Cells.Item(11, 2).Value = "extr2"
$ws.Cells.Item(12, 2).Value = "extr3"
$ws.Cells.Item(13, 2).Value = "extr4"
$ws.Cells.Item(14, 2).Value = "extr5"
$ws.Cells.Item(15, 2).Value = "extr6"
$ws.Cells.Item(16, 2).Value = "extr7"
$ws.Cells.Item(17, 2).Value = "extr8"

# --- Columns C, D, E: final numeric / boolean values for rows 8-17 ---
$values = @(
  @(8,  14, 11, $true),
  @(9,  16, 9,  $false),
  @(10, 5,  12, $true),
  @(11, 5,  9,  $true),
  @(12, 10, 11, $false),
  @(13, 7,  8,  $true),
  @(14, 9,  11, $true),
  @(15, 7,  11, $true),
  @(16, 5,  7,  $false),
  @(17, 8,  5,  $false)
)

foreach ($row in $values) {
  $r = $row[0]
  $ws.Cells.Item($r, 3).Value = $row[1]
  $ws.Cells.Item($r, 4).Value = $row[2]
  $ws.Cells.Item($r, 5).Value = $row[3]
}
